# Bridge_Intake_Form.docx — "Index Floor" field edit
#
# {Deal__r.Index_Floor__c}  ->  {Deal__r.Floor__c}%
#
# i.e. the merge field name changes (Index_Floor__c -> Floor__c) and a
# literal "%" is appended right after it, on its own run, underlined the
# same as its neighbours.

$d = $word.ActiveDocument

# --- Step 1: rename the merge field --------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "{Deal__r.Index_Floor__c}", $false, $false, $false, $false, $false,
    $true, 1, $false, "{Deal__r.Floor__c}", 2)

if (-not $found) {
    throw "Could not find '{Deal__r.Index_Floor__c}' to replace."
}

# --- Step 2: append a literal '%' right after the field -------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("{Deal__r.Floor__c}")
if (-not $found2) {
    throw "Could not find '{Deal__r.Floor__c}' after rename."
}

# Collapse to the point right after the field text, then insert "%" there.
$rng2.Collapse(0)
$rng2.InsertAfter("%")

# The inserted text inherits the surrounding (underlined) formatting, so
# nudging the underline off then back on gives "%" its own run instead of
# silently merging into a neighbouring one — matching the single-underline
# formatting of the field it follows.
$rng2.Font.Underline = 0
$rng2.Font.Underline = 1
